# Fruta / hortaliza, semanal
# Reshuffle the weekly price-report rows: for each data row, the values in
# columns D (Fecha), L (Calidad), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), Q (Unidad de
# comercializacion), R (Origen), S (Precio $/Kg) and T (Kg / unidad) are
# redistributed among the rows according to the mapping below (rows 10-12
# keep their own data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (by number) that move together as one "data block" per row.
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)

# Target row -> source row (source row's values, as they exist in the
# *original* workbook, are copied into the target row).
$rowMap = @{
    2  = 16
    3  = 7
    4  = 8
    5  = 15
    6  = 3
    7  = 4
    8  = 2
    9  = 17
    10 = 10
    11 = 11
    12 = 12
    13 = 14
    14 = 5
    15 = 6
    16 = 9
    17 = 13
}

# 1. Snapshot the original values for every row/column we might read from,
#    before any writes happen (several rows feed into each other).
$snapshot = @{}
foreach ($r in $rowMap.Keys) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value()
    }
    $snapshot[$r] = $rowVals
}

# 2. Apply the snapshot values according to the mapping.
foreach ($r in $rowMap.Keys) {
    $src = $rowMap[$r]
    $srcVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value = $srcVals[$c]
    }
}
